$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 7.786125762494509
$ws.Cells.Item(2, 4).Value = 9.548999784113093
$ws.Cells.Item(2, 5).Value = 13.75782453814284
$ws.Cells.Item(2, 6).Value = 23.71602032057022
$ws.Cells.Item(2, 7).Value = 21.44285660442219
$ws.Cells.Item(2, 8).Value = 12.374893069323
$ws.Cells.Item(2, 10).Value = 9.721307662064742
$ws.Cells.Item(2, 13).Value = 42.94676191689595
$ws.Cells.Item(2, 15).Value = 17.86342392662762

$ws.Cells.Item(3, 2).Value = 7.713342808772388
$ws.Cells.Item(3, 4).Value = 9.646280310233262
$ws.Cells.Item(3, 5).Value = 13.79054011205375
$ws.Cells.Item(3, 6).Value = 24.01827403035862
$ws.Cells.Item(3, 7).Value = 21.59406293032567
$ws.Cells.Item(3, 8).Value = 12.46511700664342
$ws.Cells.Item(3, 10).Value = 9.736926852917968
$ws.Cells.Item(3, 13).Value = 40.49120127923032
$ws.Cells.Item(3, 15).Value = 18.01486636710111

$ws.Cells.Item(4, 2).Value = 7.670124895763543
$ws.Cells.Item(4, 4).Value = 9.708943602025688
$ws.Cells.Item(4, 5).Value = 13.81720190190426
$ws.Cells.Item(4, 6).Value = 24.21556168698952
$ws.Cells.Item(4, 7).Value = 21.70388146720645
$ws.Cells.Item(4, 8).Value = 12.52408784307612
$ws.Cells.Item(4, 10).Value = 9.750067013653211
$ws.Cells.Item(4, 13).Value = 38.89726462035159
$ws.Cells.Item(4, 15).Value = 18.11553722273352

$ws.Cells.Item(5, 2).Value = 7.652902362747181
$ws.Cells.Item(5, 4).Value = 9.735217502409457
$ws.Cells.Item(5, 5).Value = 13.82970328372521
$ws.Cells.Item(5, 6).Value = 24.29887281002102
$ws.Cells.Item(5, 7).Value = 21.75281287259535
$ws.Cells.Item(5, 8).Value = 12.54901335697464
$ws.Cells.Item(5, 10).Value = 9.75630843498287
$ws.Cells.Item(5, 13).Value = 38.2262497041791
$ws.Cells.Item(5, 15).Value = 18.15847300937595

$ws.Cells.Item(6, 2).Value = 7.650066623235234
$ws.Cells.Item(6, 4).Value = 9.739624864270082
$ws.Cells.Item(6, 5).Value = 13.83187746737554
$ws.Cells.Item(6, 6).Value = 24.3128818137378
$ws.Cells.Item(6, 7).Value = 21.7611876686082
$ws.Cells.Item(6, 8).Value = 12.55320612451487
$ws.Cells.Item(6, 10).Value = 9.757398192970156
$ws.Cells.Item(6, 13).Value = 38.11353847040436
$ws.Cells.Item(6, 15).Value = 18.16571732038216

$ws.Cells.Item(7, 2).Value = 7.669891026282033
$ws.Cells.Item(7, 4).Value = 9.709294951510552
$ws.Cells.Item(7, 5).Value = 13.81736389632658
$ws.Cells.Item(7, 6).Value = 24.21667348384601
$ws.Cells.Item(7, 7).Value = 21.70452457044015
$ws.Cells.Item(7, 8).Value = 12.52442038101448
$ws.Cells.Item(7, 10).Value = 9.750147605479487
$ws.Cells.Item(7, 13).Value = 38.88830170634814
$ws.Cells.Item(7, 15).Value = 18.11610855685873

$ws.Cells.Item(8, 2).Value = 7.760736589169751
$ws.Cells.Item(8, 4).Value = 9.581933810893561
$ws.Cells.Item(8, 5).Value = 13.76773055312706
$ws.Cells.Item(8, 6).Value = 23.81779132840065
$ws.Cells.Item(8, 7).Value = 21.49141566624798
$ws.Cells.Item(8, 8).Value = 12.40525825720596
$ws.Cells.Item(8, 10).Value = 9.725952872825351
$ws.Cells.Item(8, 13).Value = 42.11808428428824
$ws.Cells.Item(8, 15).Value = 17.91403374696618

$ws.Cells.Item(9, 2).Value = 7.949683600102378
$ws.Cells.Item(9, 4).Value = 9.355420741912804
$ws.Cells.Item(9, 5).Value = 13.72327591558696
$ws.Cells.Item(9, 6).Value = 23.12975191428477
$ws.Cells.Item(9, 7).Value = 21.21209579392869
$ws.Cells.Item(9, 8).Value = 12.2001156988587
$ws.Cells.Item(9, 10).Value = 9.706929126306205
$ws.Cells.Item(9, 13).Value = 47.76124318020856
$ws.Cells.Item(9, 15).Value = 17.57968531603893

$ws.Cells.Item(10, 2).Value = 8.093883091017689
$ws.Cells.Item(10, 4).Value = 9.203136626854347
$ws.Cells.Item(10, 5).Value = 13.72376929211339
$ws.Cells.Item(10, 6).Value = 22.68355880638525
$ws.Cells.Item(10, 7).Value = 21.09643408083393
$ws.Cells.Item(10, 8).Value = 12.06705114987748
$ws.Cells.Item(10, 10).Value = 9.710610310489956
$ws.Cells.Item(10, 13).Value = 51.48082032157326
$ws.Cells.Item(10, 15).Value = 17.37311014759414

$ws.Cells.Item(11, 2).Value = 8.16038953242051
$ws.Cells.Item(11, 4).Value = 9.136923825305015
$ws.Cells.Item(11, 5).Value = 13.73136991667118
$ws.Cells.Item(11, 6).Value = 22.49389979618408
$ws.Cells.Item(11, 7).Value = 21.06426390756934
$ws.Cells.Item(11, 8).Value = 12.01040882958336
$ws.Cells.Item(11, 10).Value = 9.716182924841215
$ws.Cells.Item(11, 13).Value = 53.07983726153152
$ws.Cells.Item(11, 15).Value = 17.28791103918946

$ws.Cells.Item(12, 2).Value = 8.185682476086686
$ws.Cells.Item(12, 4).Value = 9.112291209140015
$ws.Cells.Item(12, 5).Value = 13.73532155852057
$ws.Cells.Item(12, 6).Value = 22.42403611288725
$ws.Cells.Item(12, 7).Value = 21.05509700103175
$ws.Cells.Item(12, 8).Value = 11.98952441075493
$ws.Cells.Item(12, 10).Value = 9.718858194189156
$ws.Cells.Item(12, 13).Value = 53.67194941135433
$ws.Cells.Item(12, 15).Value = 17.25693490029071

$ws.Cells.Item(13, 2).Value = 8.180230722107584
$ws.Cells.Item(13, 4).Value = 9.117576673012046
$ws.Cells.Item(13, 5).Value = 13.73442256892444
$ws.Cells.Item(13, 6).Value = 22.4389948299698
$ws.Cells.Item(13, 7).Value = 21.05693601867436
$ws.Cells.Item(13, 8).Value = 11.99399703458838
$ws.Cells.Item(13, 10).Value = 9.71825682929647
$ws.Cells.Item(13, 13).Value = 53.54502362893646
$ws.Cells.Item(13, 15).Value = 17.26354852108886

$ws.Cells.Item(14, 2).Value = 8.162468353418403
$ws.Cells.Item(14, 4).Value = 9.134888449039186
$ws.Cells.Item(14, 5).Value = 13.73167343817723
$ws.Cells.Item(14, 6).Value = 22.48811262233983
$ws.Cells.Item(14, 7).Value = 21.06344890041735
$ws.Cells.Item(14, 8).Value = 12.00867930422602
$ws.Cells.Item(14, 10).Value = 9.716391672269198
$ws.Cells.Item(14, 13).Value = 53.12881931852835
$ws.Cells.Item(14, 15).Value = 17.28533667492805

$ws.Cells.Item(15, 2).Value = 8.151601844125915
$ws.Cells.Item(15, 4).Value = 9.145549824353747
$ws.Cells.Item(15, 5).Value = 13.73012966985831
$ws.Cells.Item(15, 6).Value = 22.51845468758739
$ws.Cells.Item(15, 7).Value = 21.06783303640169
$ws.Cells.Item(15, 8).Value = 12.01774633859156
$ws.Cells.Item(15, 10).Value = 9.715322921434646
$ws.Cells.Item(15, 13).Value = 52.87213663492006
$ws.Cells.Item(15, 15).Value = 17.29885090354804

$ws.Cells.Item(16, 2).Value = 8.08955312352016
$ws.Cells.Item(16, 4).Value = 9.207525436028623
$ws.Cells.Item(16, 5).Value = 13.72342211875768
$ws.Cells.Item(16, 6).Value = 22.69622489762047
$ws.Cells.Item(16, 7).Value = 21.09895444811693
$ws.Cells.Item(16, 8).Value = 12.07083158970723
$ws.Cells.Item(16, 10).Value = 9.710324960296548
$ws.Cells.Item(16, 13).Value = 51.37444631672064
$ws.Cells.Item(16, 15).Value = 17.37885680230161

$ws.Cells.Item(17, 2).Value = 8.051705218216597
$ws.Cells.Item(17, 4).Value = 9.246329956712653
$ws.Cells.Item(17, 5).Value = 13.72120591951436
$ws.Cells.Item(17, 6).Value = 22.80872040561672
$ws.Cells.Item(17, 7).Value = 21.12333822363146
$ws.Cells.Item(17, 8).Value = 12.10439769764864
$ws.Cells.Item(17, 10).Value = 9.708260591162974
$ws.Cells.Item(17, 13).Value = 50.43181067868802
$ws.Cells.Item(17, 15).Value = 17.4302026989292

$ws.Cells.Item(18, 2).Value = 8.030023036015324
$ws.Cells.Item(18, 4).Value = 9.268937464053391
$ws.Cells.Item(18, 5).Value = 13.72062481019361
$ws.Cells.Item(18, 6).Value = 22.87467565474344
$ws.Cells.Item(18, 7).Value = 21.13928333482621
$ws.Cells.Item(18, 8).Value = 12.12406978355304
$ws.Cells.Item(18, 10).Value = 9.707439877545015
$ws.Cells.Item(18, 13).Value = 49.88087443154334
$ws.Cells.Item(18, 15).Value = 17.46055970681576

$ws.Cells.Item(19, 2).Value = 8.022697435068091
$ws.Cells.Item(19, 4).Value = 9.276641461941489
$ws.Cells.Item(19, 5).Value = 13.72054677316499
$ws.Cells.Item(19, 6).Value = 22.89722057885476
$ws.Cells.Item(19, 7).Value = 21.14500948561287
$ws.Cells.Item(19, 8).Value = 12.13079306687162
$ws.Cells.Item(19, 10).Value = 9.707224827334507
$ws.Cells.Item(19, 13).Value = 49.69283382968959
$ws.Cells.Item(19, 15).Value = 17.47097886085173

$ws.Cells.Item(20, 2).Value = 8.055725357538957
$ws.Cells.Item(20, 4).Value = 9.242169321989117
$ws.Cells.Item(20, 5).Value = 13.72136995765878
$ws.Cells.Item(20, 6).Value = 22.79661532193012
$ws.Cells.Item(20, 7).Value = 21.12054323118741
$ws.Cells.Item(20, 8).Value = 12.10078663535281
$ws.Cells.Item(20, 10).Value = 9.708442365206933
$ws.Cells.Item(20, 13).Value = 50.53306240542652
$ws.Cells.Item(20, 15).Value = 17.42465135612455

$ws.Cells.Item(21, 2).Value = 8.167682827532945
$ws.Cells.Item(21, 4).Value = 9.129791594974852
$ws.Cells.Item(21, 5).Value = 13.73245169546469
$ws.Cells.Item(21, 6).Value = 22.47363210930455
$ws.Cells.Item(21, 7).Value = 21.06145350553034
$ws.Cells.Item(21, 8).Value = 12.0043513923704
$ws.Cells.Item(21, 10).Value = 9.716924145167759
$ws.Cells.Item(21, 13).Value = 53.25143239339645
$ws.Cells.Item(21, 15).Value = 17.27890183494972

$ws.Cells.Item(22, 2).Value = 8.241475055024896
$ws.Cells.Item(22, 4).Value = 9.058915537201882
$ws.Cells.Item(22, 5).Value = 13.74595621957198
$ws.Cells.Item(22, 6).Value = 22.27396614348086
$ws.Cells.Item(22, 7).Value = 21.04044070992029
$ws.Cells.Item(22, 8).Value = 11.9446201642334
$ws.Cells.Item(22, 10).Value = 9.72576276376765
$ws.Cells.Item(22, 13).Value = 54.94993601399463
$ws.Cells.Item(22, 15).Value = 17.19115932114532

$ws.Cells.Item(23, 2).Value = 8.20204123961658
$ws.Cells.Item(23, 4).Value = 9.096508140396356
$ws.Cells.Item(23, 5).Value = 13.738171719901
$ws.Cells.Item(23, 6).Value = 22.37947228105275
$ws.Cells.Item(23, 7).Value = 21.05002138413331
$ws.Cells.Item(23, 8).Value = 11.97619650047224
$ws.Cells.Item(23, 10).Value = 9.720742526959233
$ws.Cells.Item(23, 13).Value = 54.05056224520665
$ws.Cells.Item(23, 15).Value = 17.23729318757962

$ws.Cells.Item(24, 2).Value = 8.053907610758765
$ws.Cells.Item(24, 4).Value = 9.244049415928208
$ws.Cells.Item(24, 5).Value = 13.72129363841398
$ws.Cells.Item(24, 6).Value = 22.80208404618122
$ws.Cells.Item(24, 7).Value = 21.12180085226585
$ws.Cells.Item(24, 8).Value = 12.10241803059957
$ws.Cells.Item(24, 10).Value = 9.708359045127253
$ws.Cells.Item(24, 13).Value = 50.48731453854868
$ws.Cells.Item(24, 15).Value = 17.42715851128542

$ws.Cells.Item(25, 2).Value = 7.897535158504462
$ws.Cells.Item(25, 4).Value = 9.414214409135161
$ws.Cells.Item(25, 5).Value = 13.72954279050217
$ws.Cells.Item(25, 6).Value = 23.30561402707751
$ws.Cells.Item(25, 7).Value = 21.27228556595217
$ws.Cells.Item(25, 8).Value = 12.25253135425766
$ws.Cells.Item(25, 10).Value = 9.706929126306205
$ws.Cells.Item(25, 13).Value = 46.30924321926478
$ws.Cells.Item(25, 15).Value = 17.66337174437779
